$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking strings in column D stay as text by temporarily
# forcing a Text number format, then resetting the style back to Normal so
# no extra style index is left behind on the cell (matches original formatting).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "26.668.39"
$ws.Range("E2").Value = "  +0.04%  "

# Row 3
Set-TextValue $ws.Range("D3") "1.598.90"
$ws.Range("E3").Value = "  -0.02%  "

# Row 4
$ws.Range("E4").Value = "  +0.16%  "

# Row 5
Set-TextValue $ws.Range("D5") "211.66"
$ws.Range("E5").Value = "  +0.20%  "

# Row 6
$ws.Range("E6").Value = "  +0.48%  "

# Row 7
$ws.Range("E7").Value = "  +0.16%  "

# Row 8
$ws.Range("E8").Value = "  +0.30%  "

# Row 10
Set-TextValue $ws.Range("D10") "19.53"
$ws.Range("E10").Value = "  -0.60%  "

# Row 11
$ws.Range("E11").Value = "  +0.22%  "

# Row 12
Set-TextValue $ws.Range("D12") "1.823.68"
$ws.Range("E12").Value = "  +0.00%  "

# Row 13
Set-TextValue $ws.Range("D13") "1.587.94"
$ws.Range("E13").Value = "  -0.60%  "

# Row 14
$ws.Range("E14").Value = "  -0.02%  "

# Row 16
Set-TextValue $ws.Range("D16") "65.10"
$ws.Range("E16").Value = "  +0.31%  "

# Row 17
Set-TextValue $ws.Range("D17") "26.653.13"
$ws.Range("E17").Value = "  +0.00%  "

# Row 18
Set-TextValue $ws.Range("D18") "0.0₃0739"

# Row 19
Set-TextValue $ws.Range("D19") "209.21"
$ws.Range("E19").Value = "  +0.11%  "

# Row 20
$ws.Range("E20").Value = "  +0.15%  "

# Row 21
$ws.Range("E21").Value = "  +3.93%  "

# Row 22
$ws.Range("E22").Value = "  +0.51%  "

# Row 23
Set-TextValue $ws.Range("D23") "2.36"
$ws.Range("E23").Value = "  +2.35%  "

# Row 24
Set-TextValue $ws.Range("D24") "9.00"
$ws.Range("E24").Value = "  +1.09%  "

# Row 25
Set-TextValue $ws.Range("D25") "144.33"
$ws.Range("E25").Value = "  -1.18%  "

# Row 26
$ws.Range("E26").Value = "  +0.23%  "

# Row 27
$ws.Range("E27").Value = "  -1.58%  "

# Row 28
$ws.Range("E28").Value = "  -0.62%  "

# Row 29
Set-TextValue $ws.Range("D29") "15.30"
$ws.Range("E29").Value = "  +0.06%  "

# Row 30
Set-TextValue $ws.Range("D30") "0.0514"
$ws.Range("E30").Value = "  +1.72%  "

# Row 31
$ws.Range("E31").Value = "  +0.28%  "

# Row 32
Set-TextValue $ws.Range("D32") "3.24"
$ws.Range("E32").Value = "  +0.34%  "

# Row 33
Set-TextValue $ws.Range("D33") "2.95"
$ws.Range("E33").Value = "  +1.16%  "

# Row 34
Set-TextValue $ws.Range("D34") "1.284.40"
$ws.Range("E34").Value = "  -0.85%  "

# Row 35
Set-TextValue $ws.Range("D35") "0.620"
$ws.Range("E35").Value = "  -7.15%  "

# Row 36
$ws.Range("E36").Value = "  +0.62%  "

# Row 37
$ws.Range("E37").Value = "  +0.31%  "

# Row 38
$ws.Range("E38").Value = "  -0.48%  "

# Row 39
$ws.Range("E39").Value = "  -1.08%  "

# Row 40
Set-TextValue $ws.Range("D40") "1.04"
$ws.Range("E40").Value = "  +16.18%  "

# Row 41
$ws.Range("E41").Value = "  +1.86%  "

# Row 42
$ws.Range("E42").Value = "  -0.35%  "

# Row 43
$ws.Range("E43").Value = "  -0.54%  "

# Row 44
Set-TextValue $ws.Range("D44") "63.49"
$ws.Range("E44").Value = "  -0.51%  "

# Row 45
Set-TextValue $ws.Range("D45") "1.734.51"
$ws.Range("E45").Value = "  -0.11%  "

# Row 46
Set-TextValue $ws.Range("D46") "91.10"
$ws.Range("E46").Value = "  +1.30%  "

# Row 47
$ws.Range("E47").Value = "  -2.91%  "

# Row 48
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws.Range("D48") "0.0₆0104"
$ws.Range("E48").Value = "  -0.93%  "

# Row 49
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D49") "0.102"
$ws.Range("E49").Value = "  +1.33%  "

# Row 50
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D50") "0.0509"
$ws.Range("E50").Value = "  +0.90%  "

# Row 51
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
Set-TextValue $ws.Range("D51") "1.00"
$ws.Range("E51").Value = "  +0.17%  "
